$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The municip_code column (G) held the wrong value "county03" (a county code)
# for every ward row; it should hold the municipality code "municip0301"
# instead. Fix it for all data rows (2-16).
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 7).Value = "municip0301"
}

# Column G now needs to be a bit wider to fit "municip0301" -- best-fit it.
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666

# Reflect the selection left behind after fixing/reviewing the column.
$ws.Range("G2:G16").Select() | Out-Null
